# Daily attendance processing - 2026-01-25 07:37:44
# Swap the order of "Recorded By" names in column G from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# for every row where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
